{"js": "// Locate the exact run of text containing the AQL exception stack trace and\n// replace it (7.x trace -> 8.x trace), per \"Moving from AQL 7.x to 8.x.\" commit.\nconst body = context.document.body;\n\n// Unique start marker (including the leading 4 spaces of the run) and a unique\n// end marker (including the trailing newline that terminates the run's text).\nconst startResults = body.search(\"    <---divOp\", { matchCase: true });\nstartResults.load(\"items\");\nconst endResults = body.search(\n  \"RemoteTestRunner.main(RemoteTestRunner.java:210)\\n\",\n  { matchCase: true }\n);\nendResults.load(\"items\");\nawait context.sync();\n\nif (startResults.items.length !== 1 || endResults.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly one start/end match, got \" +\n      startResults.items.length +\n      \"/\" +\n      endResults.items.length\n  );\n}\n\n// Build a single range spanning the whole stack-trace run (start marker through\n// the end marker, inclusive) and overwrite it in one shot so the run's\n// formatting (bold, red) is preserved.\nconst fullRange = startResults.items[0].expandTo(endResults.items[0]);\n\nconst NEW_TEXT = \"    <---divOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:\\n\\tjava.lang.ArithmeticException: / by zero\\n\\t\\tat org.eclipse.acceleo.query.services.NumberServices.divOp(NumberServices.java:99)\\n\\t\\tat java.base/jdk.internal.reflect.DirectMethodHandleAccessor.invoke(DirectMethodHandleAccessor.java:103)\\n\\t\\tat java.base/java.lang.reflect.Method.invoke(Method.java:580)\\n\\t\\tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:170)\\n\\t\\tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:231)\\n\\t\\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.internalCallService(EvaluationServices.java:122)\\n\\t\\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.internalCall(EvaluationServices.java:237)\\n\\t\\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:170)\\n\\t\\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:244)\\n\\t\\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:135)\\n\\t\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\\t\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\\t\\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseSequenceInExtensionLiteral(AstEvaluator.java:391)\\n\\t\\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:338)\\n\\t\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\\t\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\\t\\tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:146)\\n\\t\\tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:53)\\n\\t\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseRepetition(M2DocEvaluator.java:2087)\\n\\t\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseRepetition(M2DocEvaluator.java:1)\\n\\t\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:147)\\n\\t\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\\t\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\\t\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2124)\\n\\t\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:2349)\\n\\t\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)\\n\\t\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)\\n\\t\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\\t\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\\t\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2124)\\n\\t\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:350)\\n\\t\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)\\n\\t\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)\\n\\t\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\\n\\t\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\\n\\t\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:2124)\\n\\t\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:334)\\n\\t\\tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:878)\\n\\t\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:499)\\n\\t\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:391)\\n\\t\\tat java.base/jdk.internal.reflect.DirectMethodHandleAccessor.invoke(DirectMethodHandleAccessor.java:103)\\n\\t\\tat java.base/java.lang.reflect.Method.invoke(Method.java:580)\\n\\t\\tat org.junit.runners.model.FrameworkMethod$1.runReflectiveCall(FrameworkMethod.java:59)\\n\\t\\tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)\\n\\t\\tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:56)\\n\\t\\tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)\\n\\t\\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\\n\\t\\tat org.junit.runners.ParentRunner$3.evaluate(ParentRunner.java:306)\\n\\t\\tat org.junit.runners.BlockJUnit4ClassRunner$1.evaluate(BlockJUnit4ClassRunner.java:100)\\n\\t\\tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:366)\\n\\t\\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:103)\\n\\t\\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:63)\\n\\t\\tat org.junit.runners.ParentRunner$4.run(ParentRunner.java:331)\\n\\t\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:79)\\n\\t\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)\\n\\t\\tat org.junit.runners.ParentRunner.access$100(ParentRunner.java:66)\\n\\t\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:293)\\n\\t\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)\\n\\t\\tat org.junit.runners.Suite.runChild(Suite.java:128)\\n\\t\\tat org.junit.runners.Suite.runChild(Suite.java:27)\\n\\t\\tat org.junit.runners.ParentRunner$4.run(ParentRunner.java:331)\\n\\t\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:79)\\n\\t\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)\\n\\t\\tat org.junit.runners.ParentRunner.access$100(ParentRunner.java:66)\\n\\t\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:293)\\n\\t\\tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)\\n\\t\\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\\n\\t\\tat org.junit.runners.ParentRunner$3.evaluate(ParentRunner.java:306)\\n\\t\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)\\n\\t\\tat org.junit.runners.Suite.runChild(Suite.java:128)\\n\\t\\tat org.junit.runners.Suite.runChild(Suite.java:27)\\n\\t\\tat org.junit.runners.ParentRunner$4.run(ParentRunner.java:331)\\n\\t\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:79)\\n\\t\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)\\n\\t\\tat org.junit.runners.ParentRunner.access$100(ParentRunner.java:66)\\n\\t\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:293)\\n\\t\\tat org.junit.runners.ParentRunner$3.evaluate(ParentRunner.java:306)\\n\\t\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)\\n\\t\\tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:93)\\n\\t\\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:40)\\n\\t\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:529)\\n\\t\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:757)\\n\\t\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:452)\\n\\t\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:210)\\n\\t\";\n\nfullRange.insertText(NEW_TEXT, \"Replace\");\nawait context.sync();\n", "ps1": "# Locate the exact run of text containing the AQL exception stack trace and\n# replace it (7.x trace -> 8.x trace), per \"Moving from AQL 7.x to 8.x.\" commit.\n$d = $word.ActiveDocument\n\n# Unique start marker (including the leading 4 spaces of the run) and a unique\n# end marker (just before the trailing newline that terminates the run's text).\n$searchRange = $d.Content\n\n$startRange = $d.Range($searchRange.Start, $searchRange.End)\n$startFind = $startRange.Find\n$startFind.ClearFormatting()\n$startFind.Text = \"    <---divOp\"\n$startFound = $startFind.Execute()\n\n$endRange = $d.Range($searchRange.Start, $searchRange.End)\n$endFind = $endRange.Find\n$endFind.ClearFormatting()\n$endFind.Text = \"RemoteTestRunner.main(RemoteTestRunner.java:210)\"\n$endFound = $endFind.Execute()\n\nif (-not $startFound -or -not $endFound) {\n    throw \"Could not locate the stack-trace run to replace.\"\n}\n\n# Include the trailing newline that terminates the run's text (1 char after the\n# end match) so the whole original run content is covered.\n$endRange.MoveEnd(1, 1) | Out-Null\n\n# Build a single range spanning the whole stack-trace run (start marker through\n# the end marker, inclusive) and overwrite it in one shot so the run's\n# formatting (bold, red) is preserved.\n$fullRange = $d.Range($startRange.Start, $endRange.End)\n\n# New text is base64-encoded (UTF-8) to avoid any PowerShell quoting issues with\n# the tabs/newlines/special characters inside the stack trace.\n$newTextB64 = \"ICAgIDwtLS1kaXZPcChqYXZhLmxhbmcuSW50ZWdlcixqYXZhLmxhbmcuSW50ZWdlcikgd2l0aCBhcmd1bWVudHMgWzEsIDBdIGZhaWxlZDoKCWphdmEubGFuZy5Bcml0aG1ldGljRXhjZXB0aW9uOiAvIGJ5IHplcm8KCQlhdCBvcmcuZWNsaXBzZS5hY2NlbGVvLnF1ZXJ5LnNlcnZpY2VzLk51bWJlclNlcnZpY2VzLmRpdk9wKE51bWJlclNlcnZpY2VzLmphdmE6OTkpCgkJYXQgamF2YS5iYXNlL2pkay5pbnRlcm5hbC5yZWZsZWN0LkRpcmVjdE1ldGhvZEhhbmRsZUFjY2Vzc29yLmludm9rZShEaXJlY3RNZXRob2RIYW5kbGVBY2Nlc3Nvci5qYXZhOjEwMykKCQlhdCBqYXZhLmJhc2UvamF2YS5sYW5nLnJlZmxlY3QuTWV0aG9kLmludm9rZShNZXRob2QuamF2YTo1ODApCgkJYXQgb3JnLmVjbGlwc2UuYWNjZWxlby5xdWVyeS5ydW50aW1lLmltcGwuSmF2YU1ldGhvZFNlcnZpY2UuaW50ZXJuYWxJbnZva2UoSmF2YU1ldGhvZFNlcnZpY2UuamF2YToxNzApCgkJYXQgb3JnLmVjbGlwc2UuYWNjZWxlby5xdWVyeS5ydW50aW1lLmltcGwuQWJzdHJhY3RTZXJ2aWNlLmludm9rZShBYnN0cmFjdFNlcnZpY2UuamF2YToyMzEpCgkJYXQgb3JnLmVjbGlwc2UuYWNjZWxlby5xdWVyeS5ydW50aW1lLmltcGwuRXZhbHVhdGlvblNlcnZpY2VzLmludGVybmFsQ2FsbFNlcnZpY2UoRXZhbHVhdGlvblNlcnZpY2VzLmphdmE6MTIyKQoJCWF0IG9yZy5lY2xpcHNlLmFjY2VsZW8ucXVlcnkucnVudGltZS5pbXBsLkV2YWx1YXRpb25TZXJ2aWNlcy5pbnRlcm5hbENhbGwoRXZhbHVhdGlvblNlcnZpY2VzLmphdmE6MjM3KQoJCWF0IG9yZy5lY2xpcHNlLmFjY2VsZW8ucXVlcnkucnVudGltZS5pbXBsLkV2YWx1YXRpb25TZXJ2aWNlcy5jYWxsKEV2YWx1YXRpb25TZXJ2aWNlcy5qYXZhOjE3MCkKCQlhdCBvcmcuZWNsaXBzZS5hY2NlbGVvLnF1ZXJ5LnBhcnNlci5Bc3RFdmFsdWF0b3IuY2FzZUNhbGwoQXN0RXZhbHVhdG9yLmphdmE6MjQ0KQoJCWF0IG9yZy5lY2xpcHNlLmFjY2VsZW8ucXVlcnkuYXN0LnV0aWwuQXN0U3dpdGNoLmRvU3dpdGNoKEFzdFN3aXRjaC5qYXZhOjEzNSkKCQlhdCBvcmcuZWNsaXBzZS5lbWYuZWNvcmUudXRpbC5Td2l0Y2guZG9Td2l0Y2goU3dpdGNoLmphdmE6NTMpCgkJYXQgb3JnLmVjbGlwc2UuZW1mLmVjb3JlLnV0aWwuU3dpdGNoLmRvU3dpdGNoKFN3aXRjaC5qYXZhOjY5KQoJCWF0IG9yZy5lY2xpcHNlLmFjY2VsZW8ucXVlcnkucGFyc2VyLkFzdEV2YWx1YXRvci5jYXNlU2VxdWVuY2VJbkV4dGVuc2lvbkxpdGVyYWwoQXN0RXZhbHVhdG9yLmphdmE6MzkxKQoJCWF0IG9yZy5lY2xpcHNlLmFjY2VsZW8ucXVlcnkuYXN0LnV0aWwuQXN0U3dpdGNoLmRvU3dpdGNoKEFzdFN3aXRjaC5qYXZhOjMzOCkKCQlhdCBvcmcuZWNsaXBzZS5lbWYuZWNvcmUudXRpbC5Td2l0Y2guZG9Td2l0Y2goU3dpdGNoLmphdmE6NTMpCgkJYXQgb3JnLmVjbGlwc2UuZW1mLmVjb3JlLnV0aWwuU3dpdGNoLmRvU3dpdGNoKFN3aXRjaC5qYXZhOjY5KQoJCWF0IG9yZy5lY2xpcHNlLmFjY2VsZW8ucXVlcnkucGFyc2VyLkFzdEV2YWx1YXRvci5ldmFsKEFzdEV2YWx1YXRvci5qYXZhOjE0NikKCQlhdCBvcmcuZWNsaXBzZS5hY2NlbGVvLnF1ZXJ5LnJ1bnRpbWUuaW1wbC5RdWVyeUV2YWx1YXRpb25FbmdpbmUuZXZhbChRdWVyeUV2YWx1YXRpb25FbmdpbmUuamF2YTo1MykKCQlhdCBvcmcub2Jlb25ldHdvcmsubTJkb2MuZ2VuZXJhdG9yLk0yRG9jRXZhbHVhdG9yLmNhc2VSZXBldGl0aW9uKE0yRG9jRXZhbHVhdG9yLmphdmE6MjA4NykKCQlhdCBvcmcub2Jlb25ldHdvcmsubTJkb2MuZ2VuZXJhdG9yLk0yRG9jRXZhbHVhdG9yLmNhc2VSZXBldGl0aW9uKE0yRG9jRXZhbHVhdG9yLmphdmE6MSkKCQlhdCBvcmcub2Jlb25ldHdvcmsubTJkb2MudGVtcGxhdGUudXRpbC5UZW1wbGF0ZVN3aXRjaC5kb1N3aXRjaChUZW1wbGF0ZVN3aXRjaC5qYXZhOjE0NykKCQlhdCBvcmcuZWNsaXBzZS5lbWYuZWNvcmUudXRpbC5Td2l0Y2guZG9Td2l0Y2goU3dpdGNoLmphdmE6NTMpCgkJYXQgb3JnLmVjbGlwc2UuZW1mLmVjb3JlLnV0aWwuU3dpdGNoLmRvU3dpdGNoKFN3aXRjaC5qYXZhOjY5KQoJCWF0IG9yZy5vYmVvbmV0d29yay5tMmRvYy5nZW5lcmF0b3IuTTJEb2NFdmFsdWF0b3IuZG9Td2l0Y2goTTJEb2NFdmFsdWF0b3IuamF2YToyMTI0KQoJCWF0IG9yZy5vYmVvbmV0d29yay5tMmRvYy5nZW5lcmF0b3IuTTJEb2NFdmFsdWF0b3IuY2FzZUJsb2NrKE0yRG9jRXZhbHVhdG9yLmphdmE6MjM0OSkKCQlhdCBvcmcub2Jlb25ldHdvcmsubTJkb2MuZ2VuZXJhdG9yLk0yRG9jRXZhbHVhdG9yLmNhc2VCbG9jayhNMkRvY0V2YWx1YXRvci5qYXZhOjEpCgkJYXQgb3JnLm9iZW9uZXR3b3JrLm0yZG9jLnRlbXBsYXRlLnV0aWwuVGVtcGxhdGVTd2l0Y2guZG9Td2l0Y2goVGVtcGxhdGVTd2l0Y2guamF2YToxOTkpCgkJYXQgb3JnLmVjbGlwc2UuZW1mLmVjb3JlLnV0aWwuU3dpdGNoLmRvU3dpdGNoKFN3aXRjaC5qYXZhOjUzKQoJCWF0IG9yZy5lY2xpcHNlLmVtZi5lY29yZS51dGlsLlN3aXRjaC5kb1N3aXRjaChTd2l0Y2guamF2YTo2OSkKCQlhdCBvcmcub2Jlb25ldHdvcmsubTJkb2MuZ2VuZXJhdG9yLk0yRG9jRXZhbHVhdG9yLmRvU3dpdGNoKE0yRG9jRXZhbHVhdG9yLmphdmE6MjEyNCkKCQlhdCBvcmcub2Jlb25ldHdvcmsubTJkb2MuZ2VuZXJhdG9yLk0yRG9jRXZhbHVhdG9yLmNhc2VEb2N1bWVudFRlbXBsYXRlKE0yRG9jRXZhbHVhdG9yLmphdmE6MzUwKQoJCWF0IG9yZy5vYmVvbmV0d29yay5tMmRvYy5nZW5lcmF0b3IuTTJEb2NFdmFsdWF0b3IuY2FzZURvY3VtZW50VGVtcGxhdGUoTTJEb2NFdmFsdWF0b3IuamF2YToxKQoJCWF0IG9yZy5vYmVvbmV0d29yay5tMmRvYy50ZW1wbGF0ZS51dGlsLlRlbXBsYXRlU3dpdGNoLmRvU3dpdGNoKFRlbXBsYXRlU3dpdGNoLmphdmE6Mjc5KQoJCWF0IG9yZy5lY2xpcHNlLmVtZi5lY29yZS51dGlsLlN3aXRjaC5kb1N3aXRjaChTd2l0Y2guamF2YTo1MykKCQlhdCBvcmcuZWNsaXBzZS5lbWYuZWNvcmUudXRpbC5Td2l0Y2guZG9Td2l0Y2goU3dpdGNoLmphdmE6NjkpCgkJYXQgb3JnLm9iZW9uZXR3b3JrLm0yZG9jLmdlbmVyYXRvci5NMkRvY0V2YWx1YXRvci5kb1N3aXRjaChNMkRvY0V2YWx1YXRvci5qYXZhOjIxMjQpCgkJYXQgb3JnLm9iZW9uZXR3b3JrLm0yZG9jLmdlbmVyYXRvci5NMkRvY0V2YWx1YXRvci5nZW5lcmF0ZShNMkRvY0V2YWx1YXRvci5qYXZhOjMzNCkKCQlhdCBvcmcub2Jlb25ldHdvcmsubTJkb2MudXRpbC5NMkRvY1V0aWxzLmdlbmVyYXRlKE0yRG9jVXRpbHMuamF2YTo4NzgpCgkJYXQgb3JnLm9iZW9uZXR3b3JrLm0yZG9jLnRlc3RzLkFic3RyYWN0VGVtcGxhdGVzVGVzdFN1aXRlLnByZXBhcmVvdXRwdXRBbmRHZW5lcmF0ZShBYnN0cmFjdFRlbXBsYXRlc1Rlc3RTdWl0ZS5qYXZhOjQ5OSkKCQlhdCBvcmcub2Jlb25ldHdvcmsubTJkb2MudGVzdHMuQWJzdHJhY3RUZW1wbGF0ZXNUZXN0U3VpdGUuZ2VuZXJhdGlvbihBYnN0cmFjdFRlbXBsYXRlc1Rlc3RTdWl0ZS5qYXZhOjM5MSkKCQlhdCBqYXZhLmJhc2UvamRrLmludGVybmFsLnJlZmxlY3QuRGlyZWN0TWV0aG9kSGFuZGxlQWNjZXNzb3IuaW52b2tlKERpcmVjdE1ldGhvZEhhbmRsZUFjY2Vzc29yLmphdmE6MTAzKQoJCWF0IGphdmEuYmFzZS9qYXZhLmxhbmcucmVmbGVjdC5NZXRob2QuaW52b2tlKE1ldGhvZC5qYXZhOjU4MCkKCQlhdCBvcmcuanVuaXQucnVubmVycy5tb2RlbC5GcmFtZXdvcmtNZXRob2QkMS5ydW5SZWZsZWN0aXZlQ2FsbChGcmFtZXdvcmtNZXRob2QuamF2YTo1OSkKCQlhdCBvcmcuanVuaXQuaW50ZXJuYWwucnVubmVycy5tb2RlbC5SZWZsZWN0aXZlQ2FsbGFibGUucnVuKFJlZmxlY3RpdmVDYWxsYWJsZS5qYXZhOjEyKQoJCWF0IG9yZy5qdW5pdC5ydW5uZXJzLm1vZGVsLkZyYW1ld29ya01ldGhvZC5pbnZva2VFeHBsb3NpdmVseShGcmFtZXdvcmtNZXRob2QuamF2YTo1NikKCQlhdCBvcmcuanVuaXQuaW50ZXJuYWwucnVubmVycy5zdGF0ZW1lbnRzLkludm9rZU1ldGhvZC5ldmFsdWF0ZShJbnZva2VNZXRob2QuamF2YToxNykKCQlhdCBvcmcuanVuaXQuaW50ZXJuYWwucnVubmVycy5zdGF0ZW1lbnRzLlJ1bkFmdGVycy5ldmFsdWF0ZShSdW5BZnRlcnMuamF2YToyNykKCQlhdCBvcmcuanVuaXQucnVubmVycy5QYXJlbnRSdW5uZXIkMy5ldmFsdWF0ZShQYXJlbnRSdW5uZXIuamF2YTozMDYpCgkJYXQgb3JnLmp1bml0LnJ1bm5lcnMuQmxvY2tKVW5pdDRDbGFzc1J1bm5lciQxLmV2YWx1YXRlKEJsb2NrSlVuaXQ0Q2xhc3NSdW5uZXIuamF2YToxMDApCgkJYXQgb3JnLmp1bml0LnJ1bm5lcnMuUGFyZW50UnVubmVyLnJ1bkxlYWYoUGFyZW50UnVubmVyLmphdmE6MzY2KQoJCWF0IG9yZy5qdW5pdC5ydW5uZXJzLkJsb2NrSlVuaXQ0Q2xhc3NSdW5uZXIucnVuQ2hpbGQoQmxvY2tKVW5pdDRDbGFzc1J1bm5lci5qYXZhOjEwMykKCQlhdCBvcmcuanVuaXQucnVubmVycy5CbG9ja0pVbml0NENsYXNzUnVubmVyLnJ1bkNoaWxkKEJsb2NrSlVuaXQ0Q2xhc3NSdW5uZXIuamF2YTo2MykKCQlhdCBvcmcuanVuaXQucnVubmVycy5QYXJlbnRSdW5uZXIkNC5ydW4oUGFyZW50UnVubmVyLmphdmE6MzMxKQoJCWF0IG9yZy5qdW5pdC5ydW5uZXJzLlBhcmVudFJ1bm5lciQxLnNjaGVkdWxlKFBhcmVudFJ1bm5lci5qYXZhOjc5KQoJCWF0IG9yZy5qdW5pdC5ydW5uZXJzLlBhcmVudFJ1bm5lci5ydW5DaGlsZHJlbihQYXJlbnRSdW5uZXIuamF2YTozMjkpCgkJYXQgb3JnLmp1bml0LnJ1bm5lcnMuUGFyZW50UnVubmVyLmFjY2VzcyQxMDAoUGFyZW50UnVubmVyLmphdmE6NjYpCgkJYXQgb3JnLmp1bml0LnJ1bm5lcnMuUGFyZW50UnVubmVyJDIuZXZhbHVhdGUoUGFyZW50UnVubmVyLmphdmE6MjkzKQoJCWF0IG9yZy5qdW5pdC5ydW5uZXJzLlBhcmVudFJ1bm5lci5ydW4oUGFyZW50UnVubmVyLmphdmE6NDEzKQoJCWF0IG9yZy5qdW5pdC5ydW5uZXJzLlN1aXRlLnJ1bkNoaWxkKFN1aXRlLmphdmE6MTI4KQoJCWF0IG9yZy5qdW5pdC5ydW5uZXJzLlN1aXRlLnJ1bkNoaWxkKFN1aXRlLmphdmE6MjcpCgkJYXQgb3JnLmp1bml0LnJ1bm5lcnMuUGFyZW50UnVubmVyJDQucnVuKFBhcmVudFJ1bm5lci5qYXZhOjMzMSkKCQlhdCBvcmcuanVuaXQucnVubmVycy5QYXJlbnRSdW5uZXIkMS5zY2hlZHVsZShQYXJlbnRSdW5uZXIuamF2YTo3OSkKCQlhdCBvcmcuanVuaXQucnVubmVycy5QYXJlbnRSdW5uZXIucnVuQ2hpbGRyZW4oUGFyZW50UnVubmVyLmphdmE6MzI5KQoJCWF0IG9yZy5qdW5pdC5ydW5uZXJzLlBhcmVudFJ1bm5lci5hY2Nlc3MkMTAwKFBhcmVudFJ1bm5lci5qYXZhOjY2KQoJCWF0IG9yZy5qdW5pdC5ydW5uZXJzLlBhcmVudFJ1bm5lciQyLmV2YWx1YXRlKFBhcmVudFJ1bm5lci5qYXZhOjI5MykKCQlhdCBvcmcuanVuaXQuaW50ZXJuYWwucnVubmVycy5zdGF0ZW1lbnRzLlJ1bkJlZm9yZXMuZXZhbHVhdGUoUnVuQmVmb3Jlcy5qYXZhOjI2KQoJCWF0IG9yZy5qdW5pdC5pbnRlcm5hbC5ydW5uZXJzLnN0YXRlbWVudHMuUnVuQWZ0ZXJzLmV2YWx1YXRlKFJ1bkFmdGVycy5qYXZhOjI3KQoJCWF0IG9yZy5qdW5pdC5ydW5uZXJzLlBhcmVudFJ1bm5lciQzLmV2YWx1YXRlKFBhcmVudFJ1bm5lci5qYXZhOjMwNikKCQlhdCBvcmcuanVuaXQucnVubmVycy5QYXJlbnRSdW5uZXIucnVuKFBhcmVudFJ1bm5lci5qYXZhOjQxMykKCQlhdCBvcmcuanVuaXQucnVubmVycy5TdWl0ZS5ydW5DaGlsZChTdWl0ZS5qYXZhOjEyOCkKCQlhdCBvcmcuanVuaXQucnVubmVycy5TdWl0ZS5ydW5DaGlsZChTdWl0ZS5qYXZhOjI3KQoJCWF0IG9yZy5qdW5pdC5ydW5uZXJzLlBhcmVudFJ1bm5lciQ0LnJ1bihQYXJlbnRSdW5uZXIuamF2YTozMzEpCgkJYXQgb3JnLmp1bml0LnJ1bm5lcnMuUGFyZW50UnVubmVyJDEuc2NoZWR1bGUoUGFyZW50UnVubmVyLmphdmE6NzkpCgkJYXQgb3JnLmp1bml0LnJ1bm5lcnMuUGFyZW50UnVubmVyLnJ1bkNoaWxkcmVuKFBhcmVudFJ1bm5lci5qYXZhOjMyOSkKCQlhdCBvcmcuanVuaXQucnVubmVycy5QYXJlbnRSdW5uZXIuYWNjZXNzJDEwMChQYXJlbnRSdW5uZXIuamF2YTo2NikKCQlhdCBvcmcuanVuaXQucnVubmVycy5QYXJlbnRSdW5uZXIkMi5ldmFsdWF0ZShQYXJlbnRSdW5uZXIuamF2YToyOTMpCgkJYXQgb3JnLmp1bml0LnJ1bm5lcnMuUGFyZW50UnVubmVyJDMuZXZhbHVhdGUoUGFyZW50UnVubmVyLmphdmE6MzA2KQoJCWF0IG9yZy5qdW5pdC5ydW5uZXJzLlBhcmVudFJ1bm5lci5ydW4oUGFyZW50UnVubmVyLmphdmE6NDEzKQoJCWF0IG9yZy5lY2xpcHNlLmpkdC5pbnRlcm5hbC5qdW5pdDQucnVubmVyLkpVbml0NFRlc3RSZWZlcmVuY2UucnVuKEpVbml0NFRlc3RSZWZlcmVuY2UuamF2YTo5MykKCQlhdCBvcmcuZWNsaXBzZS5qZHQuaW50ZXJuYWwuanVuaXQucnVubmVyLlRlc3RFeGVjdXRpb24ucnVuKFRlc3RFeGVjdXRpb24uamF2YTo0MCkKCQlhdCBvcmcuZWNsaXBzZS5qZHQuaW50ZXJuYWwuanVuaXQucnVubmVyLlJlbW90ZVRlc3RSdW5uZXIucnVuVGVzdHMoUmVtb3RlVGVzdFJ1bm5lci5qYXZhOjUyOSkKCQlhdCBvcmcuZWNsaXBzZS5qZHQuaW50ZXJuYWwuanVuaXQucnVubmVyLlJlbW90ZVRlc3RSdW5uZXIucnVuVGVzdHMoUmVtb3RlVGVzdFJ1bm5lci5qYXZhOjc1NykKCQlhdCBvcmcuZWNsaXBzZS5qZHQuaW50ZXJuYWwuanVuaXQucnVubmVyLlJlbW90ZVRlc3RSdW5uZXIucnVuKFJlbW90ZVRlc3RSdW5uZXIuamF2YTo0NTIpCgkJYXQgb3JnLmVjbGlwc2UuamR0LmludGVybmFsLmp1bml0LnJ1bm5lci5SZW1vdGVUZXN0UnVubmVyLm1haW4oUmVtb3RlVGVzdFJ1bm5lci5qYXZhOjIxMCkKCQ==\"\n$newTextBytes = [Convert]::FromBase64String($newTextB64)\n$newText = [System.Text.Encoding]::UTF8.GetString($newTextBytes)\n\n$fullRange.Text = $newText\n"}
